# UC3.4.2_TC1 - update evaluations on the QuantitativeMetrics sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 5 (Compilation success): result flips to "no" and gains an explanatory note
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Called wrong method"

# Row 6 (Runtime without error): result cleared (no longer applicable)
$ws.Range("B6").Value = $null

# Row 7 (Assertion validity): result and note cleared (no longer applicable)
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null

# Row 12 (Code BLEU): refreshed score + detail breakdown
$ws.Range("B12").Value = 0.2886435712242112
$ws.Range("C12").Value = "{'codebleu': 0.28864357122421125, 'ngram_match_score': 0.08515643569758648, 'weighted_ngram_match_score': 0.0929472609639645, 'syntax_match_score': 0.6, 'dataflow_match_score': 0.3764705882352941}"

# Selection moves to B6 as captured in the saved view state
$ws.Range("B6").Select()
